$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Append 8 more rows (107-114) replicating the existing row 106 pattern
# (Username=moses, Password=bro, ID=1234, Email=m@g.c, Gender=Male, balance=0)
$source = $ws.Range("A106:F106")
for ($r = 107; $r -le 114; $r++) {
    $target = $ws.Range("A" + $r + ":F" + $r)
    $source.Copy($target)
}
